$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.860268666666666
$ws.Range("H2").Value = 14.580806
$ws.Range("I2").Value = 0.03812861294563102
$ws.Range("J2").Value = 0.03812861294563102
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 87.82520466666666
$ws.Range("N2").Value = 263.475614
$ws.Range("O2").Value = 0.4113681414249258
$ws.Range("P2").Value = 0.4113681414249258
$ws.Range("Q2").Value = 426.854090384987
$ws.Range("R2").Value = 3841.686813464884
$ws.Range("S2").Value = 0.0156848966425546
$ws.Range("T2").Value = 0.0156848966425546

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.860268666666666
$ws.Range("H3").Value = 14.580806
$ws.Range("I3").Value = 0.03812861294563102
$ws.Range("J3").Value = 0.03812861294563102
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 22.364335
$ws.Range("N3").Value = 67.09300499999999
$ws.Range("O3").Value = 0.1047532420570173
$ws.Range("P3").Value = 0.1047532420570173
$ws.Range("Q3").Value = 108.6966766513366
$ws.Range("R3").Value = 978.2700898620297
$ws.Range("S3").Value = 0.003994095821192009
$ws.Range("T3").Value = 0.003994095821192009

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.860268666666666
$ws.Range("H4").Value = 14.580806
$ws.Range("I4").Value = 0.03812861294563102
$ws.Range("J4").Value = 0.03812861294563102
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.00656633333333
$ws.Range("N4").Value = 96.019699
$ws.Range("O4").Value = 0.1499168918069617
$ws.Range("P4").Value = 0.1499168918069617
$ws.Range("Q4").Value = 155.5605114774882
$ws.Range("R4").Value = 1400.044603297394
$ws.Range("S4").Value = 0.005716123141719686
$ws.Range("T4").Value = 0.005716123141719686

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.860268666666666
$ws.Range("H5").Value = 14.580806
$ws.Range("I5").Value = 0.03812861294563102
$ws.Range("J5").Value = 0.03812861294563102
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 36.24501033333333
$ws.Range("N5").Value = 108.735031
$ws.Range("O5").Value = 0.1697695165452834
$ws.Range("P5").Value = 0.1697695165452834
$ws.Range("Q5").Value = 176.1604880461095
$ws.Range("R5").Value = 1585.444392414986
$ws.Range("S5").Value = 0.006473076186322011
$ws.Range("T5").Value = 0.006473076186322011

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.860268666666666
$ws.Range("H6").Value = 14.580806
$ws.Range("I6").Value = 0.03812861294563102
$ws.Range("J6").Value = 0.03812861294563102
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.206563333333333
$ws.Range("N6").Value = 21.61969
$ws.Range("O6").Value = 0.03375512275486358
$ws.Range("P6").Value = 0.03375512275486358
$ws.Range("Q6").Value = 35.02583396334888
$ws.Range("R6").Value = 315.23250567014
$ws.Range("S6").Value = 0.001287036010452456
$ws.Range("T6").Value = 0.001287036010452456

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4.860268666666666
$ws.Range("H7").Value = 14.580806
$ws.Range("I7").Value = 0.03812861294563102
$ws.Range("J7").Value = 0.03812861294563102
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 27.84771733333333
$ws.Range("N7").Value = 83.54315199999999
$ws.Range("O7").Value = 0.1304370854109484
$ws.Range("P7").Value = 0.1304370854109484
$ws.Range("Q7").Value = 135.3473879933902
$ws.Range("R7").Value = 1218.126491940512
$ws.Range("S7").Value = 0.004973385143390267
$ws.Range("T7").Value = 0.004973385143390265

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.24352166666667
$ws.Range("H8").Value = 42.730565
$ws.Range("I8").Value = 0.1117398567564185
$ws.Range("J8").Value = 0.1117398567564185
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 87.82520466666666
$ws.Range("N8").Value = 263.475614
$ws.Range("O8").Value = 0.4113681414249258
$ws.Range("P8").Value = 0.4113681414249258
$ws.Range("Q8").Value = 1250.940205549101
$ws.Range("R8").Value = 11258.46184994191
$ws.Range("S8").Value = 0.04596621719697533
$ws.Range("T8").Value = 0.04596621719697533

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.24352166666667
$ws.Range("H9").Value = 42.730565
$ws.Range("I9").Value = 0.1117398567564185
$ws.Range("J9").Value = 0.1117398567564185
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 22.364335
$ws.Range("N9").Value = 67.09300499999999
$ws.Range("O9").Value = 0.1047532420570173
$ws.Range("P9").Value = 0.1047532420570173
$ws.Range("Q9").Value = 318.5468901330916
$ws.Range("R9").Value = 2866.922011197824
$ws.Range("S9").Value = 0.01170511226222155
$ws.Range("T9").Value = 0.01170511226222155

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.24352166666667
$ws.Range("H10").Value = 42.730565
$ws.Range("I10").Value = 0.1117398567564185
$ws.Range("J10").Value = 0.1117398567564185
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 32.00656633333333
$ws.Range("N10").Value = 96.019699
$ws.Range("O10").Value = 0.1499168918069617
$ws.Range("P10").Value = 0.1499168918069617
$ws.Range("Q10").Value = 455.8862210444372
$ws.Range("R10").Value = 4102.975989399935
$ws.Range("S10").Value = 0.0167516920158774
$ws.Range("T10").Value = 0.0167516920158774

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 14.24352166666667
$ws.Range("H11").Value = 42.730565
$ws.Range("I11").Value = 0.1117398567564185
$ws.Range("J11").Value = 0.1117398567564185
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 36.24501033333333
$ws.Range("N11").Value = 108.735031
$ws.Range("O11").Value = 0.1697695165452834
$ws.Range("P11").Value = 0.1697695165452834
$ws.Range("Q11").Value = 516.2565899913906
$ws.Range("R11").Value = 4646.309309922515
$ws.Range("S11").Value = 0.01897002146037639
$ws.Range("T11").Value = 0.01897002146037639

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 14.24352166666667
$ws.Range("H12").Value = 42.730565
$ws.Range("I12").Value = 0.1117398567564185
$ws.Range("J12").Value = 0.1117398567564185
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 7.206563333333333
$ws.Range("N12").Value = 21.61969
$ws.Range("O12").Value = 0.03375512275486358
$ws.Range("P12").Value = 0.03375512275486358
$ws.Range("Q12").Value = 102.6468409805389
$ws.Range("R12").Value = 923.8215688248499
$ws.Range("S12").Value = 0.00377179258142378
$ws.Range("T12").Value = 0.00377179258142378

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 14.24352166666667
$ws.Range("H13").Value = 42.730565
$ws.Range("I13").Value = 0.1117398567564185
$ws.Range("J13").Value = 0.1117398567564185
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 27.84771733333333
$ws.Range("N13").Value = 83.54315199999999
$ws.Range("O13").Value = 0.1304370854109484
$ws.Range("P13").Value = 0.1304370854109484
$ws.Range("Q13").Value = 396.6495652045422
$ws.Range("R13").Value = 3569.84608684088
$ws.Range("S13").Value = 0.0145750212395441
$ws.Range("T13").Value = 0.0145750212395441

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 19.98160133333333
$ws.Range("H14").Value = 59.944804
$ws.Range("I14").Value = 0.1567548618243542
$ws.Range("J14").Value = 0.1567548618243542
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 87.82520466666666
$ws.Range("N14").Value = 263.475614
$ws.Range("O14").Value = 0.4113681414249258
$ws.Range("P14").Value = 0.4113681414249258
$ws.Range("Q14").Value = 1754.888226667739
$ws.Range("R14").Value = 15793.99404000966
$ws.Range("S14").Value = 0.06448395616800563
$ws.Range("T14").Value = 0.06448395616800563

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 19.98160133333333
$ws.Range("H15").Value = 59.944804
$ws.Range("I15").Value = 0.1567548618243542
$ws.Range("J15").Value = 0.1567548618243542
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 22.364335
$ws.Range("N15").Value = 67.09300499999999
$ws.Range("O15").Value = 0.1047532420570173
$ws.Range("P15").Value = 0.1047532420570173
$ws.Range("Q15").Value = 446.8752260551133
$ws.Range("R15").Value = 4021.877034496019
$ws.Range("S15").Value = 0.01642057998430087
$ws.Range("T15").Value = 0.01642057998430087

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 19.98160133333333
$ws.Range("H16").Value = 59.944804
$ws.Range("I16").Value = 0.1567548618243542
$ws.Range("J16").Value = 0.1567548618243542
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 32.00656633333333
$ws.Range("N16").Value = 96.019699
$ws.Range("O16").Value = 0.1499168918069617
$ws.Range("P16").Value = 0.1499168918069617
$ws.Range("Q16").Value = 639.5424485215551
$ws.Range("R16").Value = 5755.882036693996
$ws.Range("S16").Value = 0.02350020166033694
$ws.Range("T16").Value = 0.02350020166033694

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 19.98160133333333
$ws.Range("H17").Value = 59.944804
$ws.Range("I17").Value = 0.1567548618243542
$ws.Range("J17").Value = 0.1567548618243542
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 36.24501033333333
$ws.Range("N17").Value = 108.735031
$ws.Range("O17").Value = 0.1697695165452834
$ws.Range("P17").Value = 0.1697695165452834
$ws.Range("Q17").Value = 724.2333468032137
$ws.Range("R17").Value = 6518.100121228924
$ws.Range("S17").Value = 0.0266121971080433
$ws.Range("T17").Value = 0.0266121971080433

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 19.98160133333333
$ws.Range("H18").Value = 59.944804
$ws.Range("I18").Value = 0.1567548618243542
$ws.Range("J18").Value = 0.1567548618243542
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 7.206563333333333
$ws.Range("N18").Value = 21.61969
$ws.Range("O18").Value = 0.03375512275486358
$ws.Range("P18").Value = 0.03375512275486358
$ws.Range("Q18").Value = 143.9986755100844
$ws.Range("R18").Value = 1295.98807959076
$ws.Range("S18").Value = 0.005291279603302754
$ws.Range("T18").Value = 0.005291279603302754

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 19.98160133333333
$ws.Range("H19").Value = 59.944804
$ws.Range("I19").Value = 0.1567548618243542
$ws.Range("J19").Value = 0.1567548618243542
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 27.84771733333333
$ws.Range("N19").Value = 83.54315199999999
$ws.Range("O19").Value = 0.1304370854109484
$ws.Range("P19").Value = 0.1304370854109484
$ws.Range("Q19").Value = 556.4419857980231
$ws.Range("R19").Value = 5007.977872182208
$ws.Range("S19").Value = 0.0204466473003647
$ws.Range("T19").Value = 0.0204466473003647

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 23.937254
$ws.Range("H20").Value = 71.811762
$ws.Range("I20").Value = 0.1877867984967206
$ws.Range("J20").Value = 0.1877867984967206
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 87.82520466666666
$ws.Range("N20").Value = 263.475614
$ws.Range("O20").Value = 0.4113681414249258
$ws.Range("P20").Value = 0.4113681414249258
$ws.Range("Q20").Value = 2102.294231707985
$ws.Range("R20").Value = 18920.64808537187
$ws.Range("S20").Value = 0.07724950628173298
$ws.Range("T20").Value = 0.07724950628173298

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 23.937254
$ws.Range("H21").Value = 71.811762
$ws.Range("I21").Value = 0.1877867984967206
$ws.Range("J21").Value = 0.1877867984967206
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 22.364335
$ws.Range("N21").Value = 67.09300499999999
$ws.Range("O21").Value = 0.1047532420570173
$ws.Range("P21").Value = 0.1047532420570173
$ws.Range("Q21").Value = 535.3407674360899
$ws.Range("R21").Value = 4818.06690692481
$ws.Range("S21").Value = 0.0196712759580393
$ws.Range("T21").Value = 0.0196712759580393

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 23.937254
$ws.Range("H22").Value = 71.811762
$ws.Range("I22").Value = 0.1877867984967206
$ws.Range("J22").Value = 0.1877867984967206
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 32.00656633333333
$ws.Range("N22").Value = 96.019699
$ws.Range("O22").Value = 0.1499168918069617
$ws.Range("P22").Value = 0.1499168918069617
$ws.Range("Q22").Value = 766.1493079888486
$ws.Range("R22").Value = 6895.343771899638
$ws.Range("S22").Value = 0.02815241315300858
$ws.Range("T22").Value = 0.02815241315300858

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 23.937254
$ws.Range("H23").Value = 71.811762
$ws.Range("I23").Value = 0.1877867984967206
$ws.Range("J23").Value = 0.1877867984967206
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 36.24501033333333
$ws.Range("N23").Value = 108.735031
$ws.Range("O23").Value = 0.1697695165452834
$ws.Range("P23").Value = 0.1697695165452834
$ws.Range("Q23").Value = 867.6060185816247
$ws.Range("R23").Value = 7808.454167234622
$ws.Range("S23").Value = 0.03188047399437479
$ws.Range("T23").Value = 0.03188047399437479

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 23.937254
$ws.Range("H24").Value = 71.811762
$ws.Range("I24").Value = 0.1877867984967206
$ws.Range("J24").Value = 0.1877867984967206
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 7.206563333333333
$ws.Range("N24").Value = 21.61969
$ws.Range("O24").Value = 0.03375512275486358
$ws.Range("P24").Value = 0.03375512275486358
$ws.Range("Q24").Value = 172.5053369770866
$ws.Range("R24").Value = 1552.54803279378
$ws.Range("S24").Value = 0.006338766434999634
$ws.Range("T24").Value = 0.006338766434999634

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 23.937254
$ws.Range("H25").Value = 71.811762
$ws.Range("I25").Value = 0.1877867984967206
$ws.Range("J25").Value = 0.1877867984967206
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 27.84771733333333
$ws.Range("N25").Value = 83.54315199999999
$ws.Range("O25").Value = 0.1304370854109484
$ws.Range("P25").Value = 0.1304370854109484
$ws.Range("Q25").Value = 666.5978831282026
$ws.Range("R25").Value = 5999.380948153824
$ws.Range("S25").Value = 0.02449436267456529
$ws.Range("T25").Value = 0.02449436267456529

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 58.46038166666667
$ws.Range("H26").Value = 175.381145
$ws.Range("I26").Value = 0.4586193517468508
$ws.Range("J26").Value = 0.4586193517468508
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 87.82520466666666
$ws.Range("N26").Value = 263.475614
$ws.Range("O26").Value = 0.4113681414249258
$ws.Range("P26").Value = 0.4113681414249258
$ws.Range("Q26").Value = 5134.294984766448
$ws.Range("R26").Value = 46208.65486289803
$ws.Range("S26").Value = 0.1886613903496063
$ws.Range("T26").Value = 0.1886613903496063

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 58.46038166666667
$ws.Range("H27").Value = 175.381145
$ws.Range("I27").Value = 0.4586193517468508
$ws.Range("J27").Value = 0.4586193517468508
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 22.364335
$ws.Range("N27").Value = 67.09300499999999
$ws.Range("O27").Value = 0.1047532420570173
$ws.Range("P27").Value = 0.1047532420570173
$ws.Range("Q27").Value = 1307.427559821192
$ws.Range("R27").Value = 11766.84803839072
$ws.Range("S27").Value = 0.0480418639655702
$ws.Range("T27").Value = 0.0480418639655702

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 58.46038166666667
$ws.Range("H28").Value = 175.381145
$ws.Range("I28").Value = 0.4586193517468508
$ws.Range("J28").Value = 0.4586193517468508
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 32.00656633333333
$ws.Range("N28").Value = 96.019699
$ws.Range("O28").Value = 0.1499168918069617
$ws.Range("P28").Value = 0.1499168918069617
$ws.Range("Q28").Value = 1871.116083686151
$ws.Range("R28").Value = 16840.04475317536
$ws.Range("S28").Value = 0.06875478773641154
$ws.Range("T28").Value = 0.06875478773641154

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 58.46038166666667
$ws.Range("H29").Value = 175.381145
$ws.Range("I29").Value = 0.4586193517468508
$ws.Range("J29").Value = 0.4586193517468508
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 36.24501033333333
$ws.Range("N29").Value = 108.735031
$ws.Range("O29").Value = 0.1697695165452834
$ws.Range("P29").Value = 0.1697695165452834
$ws.Range("Q29").Value = 2118.897137598944
$ws.Range("R29").Value = 19070.0742383905
$ws.Range("S29").Value = 0.0778595856243741
$ws.Range("T29").Value = 0.0778595856243741

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 58.46038166666667
$ws.Range("H30").Value = 175.381145
$ws.Range("I30").Value = 0.4586193517468508
$ws.Range("J30").Value = 0.4586193517468508
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 7.206563333333333
$ws.Range("N30").Value = 21.61969
$ws.Range("O30").Value = 0.03375512275486358
$ws.Range("P30").Value = 0.03375512275486358
$ws.Range("Q30").Value = 421.2984429716722
$ws.Range("R30").Value = 3791.68598674505
$ws.Range("S30").Value = 0.01548075251597091
$ws.Range("T30").Value = 0.01548075251597091

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 58.46038166666667
$ws.Range("H31").Value = 175.381145
$ws.Range("I31").Value = 0.4586193517468508
$ws.Range("J31").Value = 0.4586193517468508
$ws.Range("K31").Value = 3
$ws.Range("M31").Value = 27.84771733333333
$ws.Range("N31").Value = 83.54315199999999
$ws.Range("O31").Value = 0.1304370854109484
$ws.Range("P31").Value = 0.1304370854109484
$ws.Range("Q31").Value = 1627.988183852116
$ws.Range("R31").Value = 14651.89365466904
$ws.Range("S31").Value = 0.05982097155491776
$ws.Range("T31").Value = 0.05982097155491774

$ws.Range("E32").Value = 3
$ws.Range("G32").Value = 5.987349666666667
$ws.Range("H32").Value = 17.962049
$ws.Range("I32").Value = 0.04697051823002506
$ws.Range("J32").Value = 0.04697051823002506
$ws.Range("K32").Value = 3
$ws.Range("M32").Value = 87.82520466666666
$ws.Range("N32").Value = 263.475614
$ws.Range("O32").Value = 0.4113681414249258
$ws.Range("P32").Value = 0.4113681414249258
$ws.Range("Q32").Value = 525.8402098858984
$ws.Range("R32").Value = 4732.561888973086
$ws.Range("S32").Value = 0.01932217478605101
$ws.Range("T32").Value = 0.01932217478605101

$ws.Range("E33").Value = 3
$ws.Range("G33").Value = 5.987349666666667
$ws.Range("H33").Value = 17.962049
$ws.Range("I33").Value = 0.04697051823002506
$ws.Range("J33").Value = 0.04697051823002506
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 22.364335
$ws.Range("N33").Value = 67.09300499999999
$ws.Range("O33").Value = 0.1047532420570173
$ws.Range("P33").Value = 0.1047532420570173
$ws.Range("Q33").Value = 133.9030937074716
$ws.Range("R33").Value = 1205.127843367245
$ws.Range("S33").Value = 0.004920314065693358
$ws.Range("T33").Value = 0.004920314065693357

$ws.Range("E34").Value = 3
$ws.Range("G34").Value = 5.987349666666667
$ws.Range("H34").Value = 17.962049
$ws.Range("I34").Value = 0.04697051823002506
$ws.Range("J34").Value = 0.04697051823002506
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 32.00656633333333
$ws.Range("N34").Value = 96.019699
$ws.Range("O34").Value = 0.1499168918069617
$ws.Range("P34").Value = 0.1499168918069617
$ws.Range("Q34").Value = 191.6345042670279
$ws.Range("R34").Value = 1724.710538403251
$ws.Range("S34").Value = 0.007041674099607591
$ws.Range("T34").Value = 0.007041674099607591

$ws.Range("E35").Value = 3
$ws.Range("G35").Value = 5.987349666666667
$ws.Range("H35").Value = 17.962049
$ws.Range("I35").Value = 0.04697051823002506
$ws.Range("J35").Value = 0.04697051823002506
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 36.24501033333333
$ws.Range("N35").Value = 108.735031
$ws.Range("O35").Value = 0.1697695165452834
$ws.Range("P35").Value = 0.1697695165452834
$ws.Range("Q35").Value = 217.0115505376132
$ws.Range("R35").Value = 1953.103954838519
$ws.Range("S35").Value = 0.007974162171792774
$ws.Range("T35").Value = 0.007974162171792774

$ws.Range("E36").Value = 3
$ws.Range("G36").Value = 5.987349666666667
$ws.Range("H36").Value = 17.962049
$ws.Range("I36").Value = 0.04697051823002506
$ws.Range("J36").Value = 0.04697051823002506
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 7.206563333333333
$ws.Range("N36").Value = 21.61969
$ws.Range("O36").Value = 0.03375512275486358
$ws.Range("P36").Value = 0.03375512275486358
$ws.Range("Q36").Value = 43.14821457164555
$ws.Range("R36").Value = 388.33393114481
$ws.Range("S36").Value = 0.001585495608714054
$ws.Range("T36").Value = 0.001585495608714054

$ws.Range("E37").Value = 3
$ws.Range("G37").Value = 5.987349666666667
$ws.Range("H37").Value = 17.962049
$ws.Range("I37").Value = 0.04697051823002506
$ws.Range("J37").Value = 0.04697051823002506
$ws.Range("K37").Value = 3
$ws.Range("M37").Value = 27.84771733333333
$ws.Range("N37").Value = 83.54315199999999
$ws.Range("O37").Value = 0.1304370854109484
$ws.Range("P37").Value = 0.1304370854109484
$ws.Range("Q37").Value = 166.7340210931609
$ws.Range("R37").Value = 1500.606189838448
$ws.Range("S37").Value = 0.006126697498166288
$ws.Range("T37").Value = 0.006126697498166286
